$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

$ws.Range("D2").Value = "48.220.59"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "2.520.76"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "322.57"
Set-TextValue $ws.Range("D6") "109.44"
$ws.Range("E6").Value = "  +1.69%  "
Set-TextValue $ws.Range("D7") "0.533"
$ws.Range("E7").Value = "  +2.36%  "
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("D9") "0.553"
Set-TextValue $ws.Range("D10") "40.52"
$ws.Range("E10").Value = "  +4.77%  "
Set-TextValue $ws.Range("D11") "20.53"
$ws.Range("E11").Value = "  +13.54%  "
Set-TextValue $ws.Range("D12") "0.0823"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E13").Value = "  +1.36%  "
Set-TextValue $ws.Range("D14") "7.28"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").Value = "2.526.08"
$ws.Range("E16").Value = "  +2.25%  "
Set-TextValue $ws.Range("D17") "0.855"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "48.051.85"
$ws.Range("E18").Value = "  +2.42%  "
Set-TextValue $ws.Range("D19") "13.31"
$ws.Range("E19").Value = "  +5.21%  "
Set-TextValue $ws.Range("D20") "6.64"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "0.0₃0947"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("E22").Value = "  -1.35%  "
Set-TextValue $ws.Range("D23") "72.01"
$ws.Range("E23").Value = "  +2.49%  "
Set-TextValue $ws.Range("D24") "264.56"
$ws.Range("E24").Value = "  +8.16%  "
Set-TextValue $ws.Range("D25") "2.58"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("E26").Value = "  -0.26%  "
Set-TextValue $ws.Range("D27") "26.05"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "2.24"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D30") "0.143"
$ws.Range("E30").Value = "  +3.83%  "
Set-TextValue $ws.Range("D31") "35.97"
$ws.Range("E31").Value = "  +3.36%  "
Set-TextValue $ws.Range("D32") "49.89"
$ws.Range("E32").Value = "  +0.94%  "
Set-TextValue $ws.Range("D33") "19.85"
$ws.Range("E33").Value = "  +0.47%  "
Set-TextValue $ws.Range("D34") "5.41"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("E37").Value = "  +2.02%  "
Set-TextValue $ws.Range("D38") "4.73"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("E40").Value = "  +0.68%  "
Set-TextValue $ws.Range("D41") "120.86"
$ws.Range("E41").Value = "  +2.11%  "
Set-TextValue $ws.Range("D42") "22.04"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("E43").Value = "  -0.79%  "
Set-TextValue $ws.Range("D44") "0.0302"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("D45").Value = "2.017.58"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("E46").Value = "  +5.93%  "
Set-TextValue $ws.Range("D48") "2.04"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  +2.21%  "
Set-TextValue $ws.Range("D51") "79.09"
$ws.Range("E51").Value = "  +3.33%  "
